# Remove envelope costs from BASoBC xlsx file
#
# The "envelope" line (row 4) on each of the three BASoBC sheets
# (urban-residential, rural-residential, commercial) is zeroed out -
# the formulas that pulled envelope-component spend from the
# Calculations sheets are replaced with literal 0s.

$wb = $excel.ActiveWorkbook

$originalActiveSheet = $wb.ActiveSheet.Name

$sheetNames = @(
    "BASoBC-urban-residential",
    "BASoBC-rural-residential",
    "BASoBC-commercial"
)

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Activate()

    $row = $ws.Range("B4:AL4")

    # Replace the "envelope" row's formulas/values with literal 0s, using
    # the same number format as the rest of the row (style index 6 in the
    # original file) so every cell ends up visually/structurally uniform.
    $row.Value = 0
    $row.NumberFormat = "0.00E+00"

    # Leave the sheet's selection on the row that was just edited.
    $row.Select()
}

# Restore the workbook's original active sheet/tab.
$wb.Worksheets.Item($originalActiveSheet).Activate()
